$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "Productos"

$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "Primer producto"
$ws.Range("C2").Value = 90
$ws.Range("D2").Value = 675
$ws.Range("E2").Value = 8
$ws.Range("F2").Value = "Camara"
$ws.Range("G2").Value = "Ayunagi"
$ws.Range("H2").Value = "Saber"
